# Auto update ESO load data
# Updates the load values in the data table (rows 2-7, columns C:AA) to the
# latest refreshed ESO load figures. Column A (Year) and B (Month) are
# unchanged; column C (Day) and the 24 hourly load columns D:AA are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,5494,5316,5240,5259,5312,5449,5895,6437,6935,7173,7138,7017,6872,6847,6791,6774,6712,7014,7448,7416,7221,6956,6723,6364),
    @(3,5936,5669,5555,5553,5594,5761,6277,6822,7340,7511,7425,7295,7153,7072,6963,6956,6997,7143,7431,7388,7164,6894,6658,6285),
    @(4,5864,5563,5433,5431,5474,5646,6177,6737,7270,7447,7358,7223,7075,6991,6879,6871,6909,7044,7311,7268,7042,6770,6533,6157),
    @(5,5514,5215,5086,5084,5127,5298,5827,6385,6915,7091,7020,6911,6792,6724,6634,6628,6657,6761,6967,6923,6693,6416,6175,5793),
    @(6,5183,4886,4759,4756,4799,4970,5496,6052,6580,6755,6700,6616,6525,6473,6403,6399,6413,6463,6562,6518,6291,6017,5779,5401),
    @(7,4983,4686,4551,4503,4527,4599,4742,5024,5609,6106,6250,6225,6158,6089,6030,6024,6053,6128,6170,6022,5776,5503,5359,5103)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colIndex = $j + 3
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
}
